$d = $word.ActiveDocument

# --- 1) Insert a new "Backend" paragraph at the very start of the body ---
$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range
$firstRange.InsertParagraphBefore()
$backendRange = $d.Range($firstRange.Start, $firstRange.Start)
$backendRange.Text = "Backend"

# --- 2) Insert "Frontend" paragraph + a following empty paragraph after the
#        paragraph that holds the long Afv6Gx... token (right before the
#        document's trailing empty paragraph) ---
$d.Content.Find.Execute("Afv6GxR1eTt5B1xU7HXVETLlhJUbiRFRypHDCyD2oy9rdk931ePbcy32cMfSCc8S9bEu6I9grHpO9zc7", $true, $false, $false, $false, $false, $true, 1, $false, "Afv6GxR1eTt5B1xU7HXVETLlhJUbiRFRypHDCyD2oy9rdk931ePbcy32cMfSCc8S9bEu6I9grHpO9zc7^pFrontend", 2)

# The replace above splits the paragraph, but the newly created run for
# "Frontend" does not inherit the "en-US" language formatting from the
# source run, so restore it explicitly.
$frontendPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$frontendPara.Range.LanguageID = "en-US"
